$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.859.99'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.894.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7842'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3144'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.34'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07201'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08094'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7657'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.482'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.60%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.918.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.39'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.169'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.868.40'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.02%  '
$ws.Range("E18").Value = '  -1.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.81'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.74%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007785'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.164'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +15.87%  '
$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.149.60'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1643'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.429'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.053'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.411'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.549'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.500'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.127'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05557'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.269'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7434'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9989'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.612'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01922'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.781'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.145.98'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +13.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.04'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4422'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.860'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8497'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.880'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.996'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.21%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.472'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.000'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +10.63%  '
